$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:D are shared-string lookups (Sending cluster, Ligand symbol,
# Receptor symbol, Target cluster); E:T are the 16 numeric metric columns.
# Each row below is an ordered array aligned to columns A..T.
$rowsData = @(
    @("ECs", "Ren1", "Atp6ap2", "ECs", 2, 0.6666666666666666, 0.2848996666666667, 0.854699, 0.3432301419625944, 0.3432301419625944, 3, 1, 7.857169666666667, 23.571509, 0.08295660471154058, 0.08295660471154058, 2.238505018976778, 20.146545170791, 0.0284732072118769, 0.02847320721187691),  # row 2
    @("ECs", "Ren1", "Atp6ap2", "FAPs", 2, 0.6666666666666666, 0.2848996666666667, 0.854699, 0.3432301419625944, 0.3432301419625944, 3, 1, 24.33072566666667, 72.992177, 0.2568856823898633, 0.2568856823898633, 6.931815632191444, 62.386340689723, 0.0881709092348307, 0.08817090923483073),  # row 3
    @("ECs", "Ren1", "Atp6ap2", "M2", 2, 0.6666666666666666, 0.2848996666666667, 0.854699, 0.3432301419625944, 0.3432301419625944, 3, 1, 53.717676, 161.153028, 0.5671553756640626, 0.5671553756640626, 15.304147986508, 137.737331878572, 0.1946648201040247, 0.1946648201040248),  # row 4
    @("ECs", "Ren1", "Atp6ap2", "sCs", 2, 0.6666666666666666, 0.2848996666666667, 0.854699, 0.3432301419625944, 0.3432301419625944, 3, 1, 8.808643333333332, 26.42593, 0.0930023372345335, 0.0930023372345335, 2.509579549452222, 22.58621594507, 0.03192120541186201, 0.03192120541186201),  # row 5
    @("FAPs", "Ren1", "Atp6ap2", "ECs", 3, 1, 0.5451546666666667, 1.635464, 0.6567698580374056, 0.6567698580374056, 3, 1, 7.857169666666667, 23.571509, 0.08295660471154058, 0.08295660471154058, 4.283372710575112, 38.550354395176, 0.05448339749966368, 0.05448339749966368),  # row 6
    @("FAPs", "Ren1", "Atp6ap2", "FAPs", 3, 1, 0.5451546666666667, 1.635464, 0.6567698580374056, 0.6567698580374056, 3, 1, 24.33072566666667, 72.992177, 0.2568856823898633, 0.2568856823898633, 13.26400864056978, 119.376077765128, 0.1687147731550326, 0.1687147731550326),  # row 7
    @("FAPs", "Ren1", "Atp6ap2", "M2", 3, 1, 0.5451546666666667, 1.635464, 0.6567698580374056, 0.6567698580374056, 3, 1, 53.717676, 161.153028, 0.5671553756640626, 0.5671553756640626, 29.284441753888, 263.559975784992, 0.3724905555600379, 0.3724905555600379),  # row 8
    @("FAPs", "Ren1", "Atp6ap2", "sCs", 3, 1, 0.5451546666666667, 1.635464, 0.6567698580374056, 0.6567698580374056, 3, 1, 8.808643333333332, 26.42593, 0.0930023372345335, 0.0930023372345335, 4.802073020168889, 43.21865718151999, 0.06108113182267148, 0.06108113182267148)  # row 9
)

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowNum = $startRow + $i
    $values = $rowsData[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
    }
}
